# Actualización automática 2025-07-09 09:35:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("D43").Value = 570.24
$wsGrupo.Range("M43").Value = 1430.71
$wsGrupo.Range("D55").Value = "2 de 53"
$wsGrupo.Range("M55").Value = "7 de 53"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F43").Value = 2027.05
$wsMensual.Range("F55").Value = 11178.3

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column F width 26 -> 25 (stored OOXML width units); ColumnWidth property
# uses a slightly different scale, so 24.1666... maps to a stored width of 25.
$wsCumpl.Columns.Item(6).ColumnWidth = 24.166666666666668

# Row 3 (240X80 PORCELANATO)
$wsCumpl.Range("D3").Value = 660
$wsCumpl.Range("E3").Value = 13068
$wsCumpl.Range("F3").Value = 0.04807692307692308

# Row 16 (PORCELANATO)
$wsCumpl.Range("D16").Value = 5600.78
$wsCumpl.Range("E16").Value = 49120.45
$wsCumpl.Range("F16").Value = 0.1023511350165192

# Row 19 (TOTAL)
$wsCumpl.Range("D19").Value = 11969.41
$wsCumpl.Range("E19").Value = 93243.45999999999
$wsCumpl.Range("F19").Value = 0.1137637439222027
